# Update the cryptos list (prices / 1h-volume percentages) and fix the
# row order for InjectiveProtocol / Stellar (their data was swapped).
#
# Columns D (Price) hold values that are sometimes plain decimal numbers
# (e.g. "2.70", "0.113"). Those look numeric to Excel, which would
# silently coerce them to real numbers (stripping formatting such as
# trailing zeros). To keep them as literal text - matching the original
# inline-string cells - we force the cell to Text format ("@") before
# assigning the value, then restore the cell's style afterwards so no
# stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.072.30'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '3.121.83'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.115.87'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.146'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.33'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.19%  '
$ws.Range('D15').Value = '3.631.00'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('D17').Value = '63.112.29'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').Value = '3.124.17'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '474.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.696'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.80%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.67%  '
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.36%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.107'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -11.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('D38').Value = '0.0₃0713'
$ws.Range('E38').Value = '  -4.70%  '
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '421.42'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.116'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -13.18%  '
$ws.Range('D44').Value = '2.863.45'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.256'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.46%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.113'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.68'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.82%  '
